# Atualização de bases das ligas, do dia: 15-06-2024 às 21:10
#
# For each listed row pair, the data columns B:AD (everything except the
# leading index column A) were swapped between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(130, 131),
    @(134, 135),
    @(137, 138),
    @(143, 144),
    @(236, 237)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($r1, 2), $ws.Cells.Item($r1, 30))
    $range2 = $ws.Range($ws.Cells.Item($r2, 2), $ws.Cells.Item($r2, 30))

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
